# Apply the changes described by the commit:
#   1. Re-style the table on slide 5 to use the built-in table style
#      {D4E3EC5E-FB16-4531-9C07-EC3AA9F21334} (was
#      {9E9F7E03-4620-499C-BAE4-F56842FF8212}).
#   2. Swap the deck's colour theme from the "Integral" / "Red Violet"
#      palette to the stock "Office Theme" / "Office" palette (the font
#      scheme and format scheme - fills/lines/effects - are identical
#      between the two themes, only the 12 theme colours differ).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{D4E3EC5E-FB16-4531-9C07-EC3AA9F21334}")
    }
}

# --- 2. Theme colours -------------------------------------------------
# Office theme colour scheme, in the standard PowerPoint
# ThemeColorScheme order: dk1, lt1, dk2, lt2, accent1..accent6, hlink,
# folHlink. Values are decimal OLE_COLOR (0x00BBGGRR) equivalents of the
# stock "Office" palette's RRGGBB hex codes, matching what
# ThemeColorScheme.Item(n).RGB reads/writes.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
